$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 205
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
